$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 639.5294
$ws.Range("I28").Value = 591.4666999999999
$ws.Range("K28").Value = 591.4666999999999
$ws.Range("M28").Value = -106.4666999999999
$ws.Range("H86").Value = 4794.385
$ws.Range("I86").Value = 4464.4443
$ws.Range("J86").Value = 5536.75
$ws.Range("K86").Value = 4464.4443
$ws.Range("L86").Value = 5536.75
$ws.Range("M86").Value = -3341.4443
$ws.Range("N86").Value = -7782.75
$ws.Range("H89").Value = 4794.385
$ws.Range("I89").Value = 4464.4443
$ws.Range("J89").Value = 5536.75
$ws.Range("K89").Value = 22322.2215
$ws.Range("L89").Value = 27683.75
$ws.Range("M89").Value = -16706.2215
$ws.Range("N89").Value = -38915.75
$ws.Range("H92").Value = 219.2
$ws.Range("I92").Value = 195.04167
$ws.Range("K92").Value = 195.04167
$ws.Range("M92").Value = 1052.95833
$ws.Range("H112").Value = 2865.5264
$ws.Range("I112").Value = 2614.6667
$ws.Range("J112").Value = 2981.3076
$ws.Range("K112").Value = 7844.000100000001
$ws.Range("L112").Value = 8943.9228
$ws.Range("M112").Value = -6736.000100000001
$ws.Range("N112").Value = -11159.9228
$ws.Range("H135").Value = 1382
$ws.Range("I135").Value = 888.2593000000001
$ws.Range("J135").Value = 4048.2
$ws.Range("K135").Value = 7994.3337
$ws.Range("L135").Value = 36433.8
$ws.Range("M135").Value = -5459.3337
$ws.Range("N135").Value = -41503.8
$ws.Range("H137").Value = 4648.2334
$ws.Range("I137").Value = 2293.4146
$ws.Range("J137").Value = 9729.684999999999
$ws.Range("K137").Value = 6880.2438
$ws.Range("L137").Value = 29189.055
$ws.Range("M137").Value = -4330.2438
$ws.Range("N137").Value = -34289.055

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 128404.086
$ws.Range("I32").Value = 149327.3
$ws.Range("K32").Value = 149327.3
$ws.Range("M32").Value = -149040.3
$ws.Range("H122").Value = 83335640
$ws.Range("I122").Value = 125002010
$ws.Range("K122").Value = 375006030
$ws.Range("M122").Value = -375003580

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1925
$ws.Range("I37").Value = 1925
$ws.Range("K37").Value = 1925
$ws.Range("M37").Value = -1788
$ws.Range("H99").Value = 7357.5
$ws.Range("I99").Value = 7357.5
$ws.Range("K99").Value = 7357.5
$ws.Range("M99").Value = -5859.5
$ws.Range("H107").Value = 9616680
$ws.Range("I107").Value = 11364928
$ws.Range("J107").Value = 1316.5
$ws.Range("K107").Value = 11364928
$ws.Range("L107").Value = 1316.5
$ws.Range("M107").Value = -11363008
$ws.Range("N107").Value = -5156.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1694.5
$ws.Range("I7").Value = 34
$ws.Range("K7").Value = 34
$ws.Range("M7").Value = 79
$ws.Range("H16").Value = 43747.625
$ws.Range("I16").Value = 2034.125
$ws.Range("J16").Value = 127174.625
$ws.Range("K16").Value = 2034.125
$ws.Range("L16").Value = 127174.625
$ws.Range("M16").Value = -1747.125
$ws.Range("N16").Value = -127748.625
$ws.Range("H31").Value = 3143.76
$ws.Range("I31").Value = 3057.342
$ws.Range("K31").Value = 3057.342
$ws.Range("M31").Value = -2762.342
$ws.Range("H34").Value = 3143.76
$ws.Range("I34").Value = 3057.342
$ws.Range("K34").Value = 3057.342
$ws.Range("M34").Value = -2855.342
$ws.Range("H44").Value = 2832.6667
$ws.Range("J44").Value = 3999
$ws.Range("L44").Value = 3999
$ws.Range("N44").Value = -4883
$ws.Range("H50").Value = 23799.8
$ws.Range("J50").Value = 23799.8
$ws.Range("L50").Value = 23799.8
$ws.Range("N50").Value = -25049.8
$ws.Range("H54").Value = 36496.5
$ws.Range("J54").Value = 36496.5
$ws.Range("L54").Value = 36496.5
$ws.Range("N54").Value = -37812.5
$ws.Range("H74").Value = 40000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 40000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H107").Value = 605.36365
$ws.Range("I107").Value = 573.2222
$ws.Range("K107").Value = 573.2222
$ws.Range("M107").Value = 1346.7778
$ws.Range("H113").Value = 43747.625
$ws.Range("I113").Value = 2034.125
$ws.Range("J113").Value = 127174.625
$ws.Range("K113").Value = 2034.125
$ws.Range("L113").Value = 127174.625
$ws.Range("M113").Value = 135.875
$ws.Range("N113").Value = -131514.625
$ws.Range("H132").Value = 1948.262
$ws.Range("I132").Value = 1742.919
$ws.Range("K132").Value = 5228.757000000001
$ws.Range("M132").Value = -2698.757000000001
$ws.Range("H134").Value = 1752
$ws.Range("I134").Value = 1560.9412
$ws.Range("K134").Value = 4682.8236
$ws.Range("M134").Value = -2147.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1478.25
$ws.Range("J92").Value = 1478.25
$ws.Range("L92").Value = 4434.75
$ws.Range("N92").Value = -6930.75
$ws.Range("H122").Value = 768473.7
$ws.Range("I122").Value = 1466549.1
$ws.Range("J122").Value = 590.7
$ws.Range("K122").Value = 13198941.9
$ws.Range("L122").Value = 5316.3
$ws.Range("M122").Value = -13196491.9
$ws.Range("N122").Value = -10216.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 546.54285
$ws.Range("I97").Value = 555.4545000000001
$ws.Range("J97").Value = 531.46155
$ws.Range("K97").Value = 555.4545000000001
$ws.Range("L97").Value = 531.46155
$ws.Range("M97").Value = -59.45450000000005
$ws.Range("N97").Value = -1523.46155
$ws.Range("H122").Value = 3257.6572
$ws.Range("I122").Value = 3150.7693
$ws.Range("K122").Value = 9452.3079
$ws.Range("M122").Value = -7002.3079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2702.6785
$ws.Range("I22").Value = 1570.2142
$ws.Range("J22").Value = 3835.1428
$ws.Range("K22").Value = 1570.2142
$ws.Range("L22").Value = 3835.1428
$ws.Range("M22").Value = -1275.2142
$ws.Range("N22").Value = -4425.1428
$ws.Range("H27").Value = 2702.6785
$ws.Range("I27").Value = 1570.2142
$ws.Range("J27").Value = 3835.1428
$ws.Range("K27").Value = 1570.2142
$ws.Range("L27").Value = 3835.1428
$ws.Range("M27").Value = -1463.2142
$ws.Range("N27").Value = -4049.1428
$ws.Range("H40").Value = 3282.7368
$ws.Range("I40").Value = 3323.1667
$ws.Range("J40").Value = 3213.4285
$ws.Range("K40").Value = 3323.1667
$ws.Range("L40").Value = 3213.4285
$ws.Range("M40").Value = -3187.1667
$ws.Range("N40").Value = -3485.4285
$ws.Range("H46").Value = 3851.6487
$ws.Range("I46").Value = 1058.8334
$ws.Range("K46").Value = 1058.8334
$ws.Range("M46").Value = -870.8334
$ws.Range("H61").Value = 4616.9165
$ws.Range("I61").Value = 4600.2607
$ws.Range("K61").Value = 4600.2607
$ws.Range("M61").Value = -4398.2607
$ws.Range("H113").Value = 4616.9165
$ws.Range("I113").Value = 4600.2607
$ws.Range("K113").Value = 4600.2607
$ws.Range("M113").Value = -2430.2607
$ws.Range("H132").Value = 3235.4666
$ws.Range("I132").Value = 2646.6086
$ws.Range("K132").Value = 7939.825800000001
$ws.Range("M132").Value = -5409.825800000001
$ws.Range("H136").Value = 3297.9048
$ws.Range("I136").Value = 1962.8
$ws.Range("J136").Value = 30000
$ws.Range("K136").Value = 5888.4
$ws.Range("L136").Value = 90000
$ws.Range("M136").Value = -3338.4
$ws.Range("N136").Value = -95100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1226.3334
$ws.Range("I100").Value = 1271.8
$ws.Range("K100").Value = 2543.6
$ws.Range("M100").Value = -2002.6
$ws.Range("H107").Value = 2615.5334
$ws.Range("I107").Value = 1516.6428
$ws.Range("J107").Value = 18000
$ws.Range("K107").Value = 4549.928400000001
$ws.Range("L107").Value = 54000
$ws.Range("M107").Value = -2629.928400000001
$ws.Range("N107").Value = -57840
$ws.Range("H136").Value = 739.44116
$ws.Range("I136").Value = 739.44116
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2218.32348
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 331.67652
$ws.Range("N136").ClearContents()
